# Update countries & provincias Spain
# - Swap order of "Brunei" / "Gibraltar" rows (Gibraltar now listed before Brunei,
#   with Gibraltar's daily numbers refreshed and Brunei's numbers carried over).
# - Refresh the daily COVID figures (B..H) for a handful of other countries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country name swap: row 138 becomes Gibraltar, row 139 becomes Brunei ---
$ws.Range("A138").Value = "Gibraltar"
$ws.Range("A139").Value = "Brunei"

# --- Row 4 : Estados Unidos ---
$ws.Range("B4").Value = 987020
$ws.Range("C4").Value = 26369
$ws.Range("E4").Value = 812836
$ws.Range("G4").Value = 1151
$ws.Range("H4").Value = 55407

# --- Row 8 : Alemania ---
$ws.Range("B8").Value = 157770
$ws.Range("C8").Value = 1257
$ws.Range("E8").Value = 39794
$ws.Range("G8").Value = 99
$ws.Range("H8").Value = 5976

# --- Row 15 : Canada ---
$ws.Range("B15").Value = 46895
$ws.Range("C15").Value = 1541
$ws.Range("D15").Value = 17321
$ws.Range("E15").Value = 27014

# --- Row 44 : Noruega ---
$ws.Range("B44").Value = 7527
$ws.Range("C44").Value = 34
$ws.Range("E44").Value = 7294

# --- Row 45 : Chequia ---
$ws.Range("B45").Value = 7404
$ws.Range("C45").Value = 52
$ws.Range("E45").Value = 4639

# --- Row 56 : Argentina ---
$ws.Range("B56").Value = 3892
$ws.Range("C56").Value = 112
$ws.Range("E56").Value = 2593
$ws.Range("G56").Value = 7
$ws.Range("H56").Value = 192

# --- Row 62 : Kazajistan ---
$ws.Range("B62").Value = 2717
$ws.Range("C62").Value = 116
$ws.Range("D62").Value = 682
$ws.Range("E62").Value = 2010

# --- Row 85 : Nigeria ---
$ws.Range("B85").Value = 1273
$ws.Range("C85").Value = 91
$ws.Range("D85").Value = 239
$ws.Range("E85").Value = 994
$ws.Range("G85").Value = 5
$ws.Range("H85").Value = 40

# --- Row 87 : Hong Kong ---
$ws.Range("D87").Value = 772
$ws.Range("E87").Value = 262

# --- Row 120 : Venezuela ---
$ws.Range("B120").Value = 325
$ws.Range("C120").Value = 2
$ws.Range("D120").Value = 137
$ws.Range("E120").Value = 178

# --- Row 138 : now Gibraltar (numbers refreshed) ---
$ws.Range("B138").Value = 141
$ws.Range("C138").Value = 5
$ws.Range("D138").Value = 131
$ws.Range("E138").Value = 10
$ws.Range("F138").Value = 0
$ws.Range("H138").Value = 0

# --- Row 139 : now Brunei (numbers carried over from previous Brunei row) ---
$ws.Range("B139").Value = 138
$ws.Range("D139").Value = 123
$ws.Range("E139").Value = 14
$ws.Range("F139").Value = 2
$ws.Range("H139").Value = 1
